$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Listado")

# Update D2 / D3 formulas: change the TEXT() format mask from "00-00000000-0" to "00000000000"
$ws.Range("D2").Formula = '=TEXT(C2,"00000000000")'
$ws.Range("D3").Formula = '=TEXT(C3,"00000000000")'

# Replace D4:D69 (currently text "30-00000000-0") with the plain numeric value 30000000000
$ws.Range("D4:D69").Value = 30000000000

# Move the active cell selection to D1
$ws.Range("D1").Select()
